# going to change my_camera_changeSize to run altered size on start-up
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 79: add DEC2BIN example (A79/B79) ---
$ws.Range("A79").Value = 2060
$ws.Range("B79").Formula = "=DEC2BIN(A79,16)"

# --- Row 80: add DEC2BIN example (A80/B80) ---
$ws.Range("A80").Value = 984
$ws.Range("B80").Formula = "=DEC2BIN(A80,32)"

# --- Row 84: HEX2BIN example with "0B" ---
$ws.Range("A84").Value = "0B"
$ws.Range("B84").Formula = "=HEX2BIN(A84)"

# --- Row 85: HEX2BIN example with "1C" ---
$ws.Range("A85").Value = "1C"
$ws.Range("B85").Formula = "=HEX2BIN(A85)"

# --- Row 87: new size calc block ---
$ws.Range("D87").Value = 2560
$ws.Range("E87").Value = 1920
$ws.Range("F87").Value = 320
$ws.Range("G87").Value = 2543
$ws.Range("H87").Value = 31
$ws.Range("I87").Value = 2684

# --- Row 88: new size calc block ---
$ws.Range("D88").Value = 1440
$ws.Range("E88").Value = 1920
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 1951
$ws.Range("H88").Value = 16
$ws.Range("I88").Value = 1968

# --- Row 92: difference check ---
$ws.Range("D92").Value = 2592
$ws.Range("E92").Value = 2684
$ws.Range("F92").Formula = "=E92-D92"

# --- Row 93: difference check ---
$ws.Range("D93").Value = 1944
$ws.Range("E93").Value = 1968
$ws.Range("F93").Formula = "=E93-D93"
$ws.Range("G93").Formula = "=F93/2"

# --- update view: scroll + selection near the new content ---
$ws.Activate()
$ws.Range("A72").Select()
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E80").Select()

# --- window position/size (matches the author's window move on save) ---
$win = $excel.ActiveWindow
$win.Left = 16365
$win.Top = 13920
$win.Width = 21675
$win.Height = 13830
